$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 739.0952
$ws.Range("I19").Value = 660.3333
$ws.Range("J19").Value = 752.2222
$ws.Range("K19").Value = 660.3333
$ws.Range("L19").Value = 752.2222
$ws.Range("M19").Value = -485.3333
$ws.Range("N19").Value = -1102.2222
$ws.Range("H43").Value = 1975.5
$ws.Range("I43").Value = 1450
$ws.Range("J43").Value = 2501
$ws.Range("K43").Value = 1450
$ws.Range("L43").Value = 2501
$ws.Range("M43").Value = -1381
$ws.Range("N43").Value = -2639
$ws.Range("H55").Value = 409
$ws.Range("I55").Value = 398.8889
$ws.Range("J55").Value = 500
$ws.Range("K55").Value = 398.8889
$ws.Range("L55").Value = 500
$ws.Range("M55").Value = -184.8889
$ws.Range("N55").Value = -928
$ws.Range("H103").Value = 1042.8334
$ws.Range("I103").Value = 700
$ws.Range("J103").Value = 1140.7858
$ws.Range("K103").Value = 2100
$ws.Range("L103").Value = 3422.3574
$ws.Range("M103").Value = -1514
$ws.Range("N103").Value = -4594.357400000001
$ws.Range("H107").Value = 737.913
$ws.Range("I107").Value = 550.58826
$ws.Range("J107").Value = 1268.6666
$ws.Range("K107").Value = 550.58826
$ws.Range("L107").Value = 1268.6666
$ws.Range("M107").Value = 1369.41174
$ws.Range("N107").Value = -5108.6666
$ws.Range("H111").Value = 3956.875
$ws.Range("I111").Value = 5097.25
$ws.Range("J111").Value = 2816.5
$ws.Range("K111").Value = 15291.75
$ws.Range("L111").Value = 8449.5
$ws.Range("M111").Value = -12224.75
$ws.Range("N111").Value = -14583.5
$ws.Range("H112").Value = 1287.4286
$ws.Range("J112").Value = 1392.258
$ws.Range("L112").Value = 4176.774
$ws.Range("N112").Value = -6392.774
$ws.Range("H113").Value = 2636.1738
$ws.Range("I113").Value = 1770.4
$ws.Range("J113").Value = 3302.1538
$ws.Range("K113").Value = 1770.4
$ws.Range("L113").Value = 3302.1538
$ws.Range("M113").Value = 1483.6
$ws.Range("N113").Value = -9810.1538
$ws.Range("H125").Value = 1286.1428
$ws.Range("I125").Value = 965.5
$ws.Range("J125").Value = 1414.4
$ws.Range("K125").Value = 8689.5
$ws.Range("L125").Value = 12729.6
$ws.Range("M125").Value = -6229.5
$ws.Range("N125").Value = -17649.6
$ws.Range("H135").Value = 68182456
$ws.Range("I135").Value = 27778290
$ws.Range("J135").Value = 250001220
$ws.Range("K135").Value = 250004610
$ws.Range("L135").Value = 2250010980
$ws.Range("M135").Value = -250002075
$ws.Range("N135").Value = -2250016050
$ws.Range("H138").Value = 1498.72
$ws.Range("I138").Value = 730.65216
$ws.Range("J138").Value = 1728.1428
$ws.Range("K138").Value = 2191.95648
$ws.Range("L138").Value = 5184.428400000001
$ws.Range("M138").Value = 2948.04352
$ws.Range("N138").Value = -15464.4284

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1568
$ws.Range("I2").Value = 1506.2858
$ws.Range("K2").Value = 1506.2858
$ws.Range("M2").Value = -1393.2858
$ws.Range("H45").Value = 1855.8422
$ws.Range("I45").Value = 1590.2333
$ws.Range("J45").Value = 2851.875
$ws.Range("K45").Value = 1590.2333
$ws.Range("L45").Value = 2851.875
$ws.Range("M45").Value = -1213.2333
$ws.Range("N45").Value = -3605.875
$ws.Range("H61").Value = 6113.143
$ws.Range("I61").Value = 4396.619
$ws.Range("K61").Value = 4396.619
$ws.Range("M61").Value = -4184.619
$ws.Range("H74").Value = 6509.5625
$ws.Range("I74").Value = 2988.6155
$ws.Range("J74").Value = 21767
$ws.Range("K74").Value = 2988.6155
$ws.Range("L74").Value = 21767
$ws.Range("M74").Value = -2114.6155
$ws.Range("N74").Value = -23515
$ws.Range("H77").Value = 6509.5625
$ws.Range("I77").Value = 2988.6155
$ws.Range("J77").Value = 21767
$ws.Range("K77").Value = 14943.0775
$ws.Range("L77").Value = 108835
$ws.Range("M77").Value = -10575.0775
$ws.Range("N77").Value = -117571
$ws.Range("H116").Value = 1568
$ws.Range("I116").Value = 1506.2858
$ws.Range("K116").Value = 1506.2858
$ws.Range("M116").Value = 787.7141999999999
$ws.Range("H136").Value = 6113.143
$ws.Range("I136").Value = 4396.619
$ws.Range("K136").Value = 13189.857
$ws.Range("M136").Value = -10639.857

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1568
$ws.Range("I3").Value = 1506.2858
$ws.Range("K3").Value = 1506.2858
$ws.Range("M3").Value = -1392.2858
$ws.Range("H20").Value = 817.94116
$ws.Range("I20").Value = 785
$ws.Range("J20").Value = 925
$ws.Range("K20").Value = 785
$ws.Range("L20").Value = 925
$ws.Range("M20").Value = -538
$ws.Range("N20").Value = -1419
$ws.Range("H26").Value = 13396.8
$ws.Range("I26").Value = 13396.8
$ws.Range("K26").Value = 13396.8
$ws.Range("M26").Value = -13104.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2845.9
$ws.Range("I31").Value = 1974.52
$ws.Range("J31").Value = 7202.8
$ws.Range("K31").Value = 1974.52
$ws.Range("L31").Value = 7202.8
$ws.Range("M31").Value = -1679.52
$ws.Range("N31").Value = -7792.8
$ws.Range("H34").Value = 2845.9
$ws.Range("I34").Value = 1974.52
$ws.Range("J34").Value = 7202.8
$ws.Range("K34").Value = 1974.52
$ws.Range("L34").Value = 7202.8
$ws.Range("M34").Value = -1772.52
$ws.Range("N34").Value = -7606.8
$ws.Range("H58").Value = 2333766.5
$ws.Range("I58").Value = 3789474.2
$ws.Range("J58").Value = 4633.8667
$ws.Range("K58").Value = 3789474.2
$ws.Range("L58").Value = 4633.8667
$ws.Range("M58").Value = -3789271.2
$ws.Range("N58").Value = -5039.8667
$ws.Range("H94").Value = 1397.2778
$ws.Range("I94").Value = 1635.1428
$ws.Range("J94").Value = 1245.909
$ws.Range("K94").Value = 1635.1428
$ws.Range("L94").Value = 1245.909
$ws.Range("M94").Value = -1184.1428
$ws.Range("N94").Value = -2147.909
$ws.Range("H132").Value = 2379.3333
$ws.Range("I132").Value = 1791.238
$ws.Range("J132").Value = 4437.6665
$ws.Range("K132").Value = 5373.714
$ws.Range("L132").Value = 13312.9995
$ws.Range("M132").Value = -2843.714
$ws.Range("N132").Value = -18372.9995
$ws.Range("H136").Value = 2333766.5
$ws.Range("I136").Value = 3789474.2
$ws.Range("J136").Value = 4633.8667
$ws.Range("K136").Value = 11368422.6
$ws.Range("L136").Value = 13901.6001
$ws.Range("M136").Value = -11365872.6
$ws.Range("N136").Value = -19001.6001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 23.84
$ws.Range("I2").Value = 18.285715
$ws.Range("J2").Value = 26
$ws.Range("K2").Value = 109.71429
$ws.Range("L2").Value = 156
$ws.Range("M2").Value = 3.285709999999995
$ws.Range("N2").Value = -382
$ws.Range("H34").Value = 1302.7307
$ws.Range("I34").Value = 205.55556
$ws.Range("J34").Value = 1883.5883
$ws.Range("K34").Value = 616.66668
$ws.Range("L34").Value = 5650.7649
$ws.Range("M34").Value = -532.66668
$ws.Range("N34").Value = -5818.7649
$ws.Range("H38").Value = 72.13333
$ws.Range("I38").Value = 33.125
$ws.Range("J38").Value = 116.71429
$ws.Range("K38").Value = 99.375
$ws.Range("L38").Value = 350.14287
$ws.Range("M38").Value = 247.625
$ws.Range("N38").Value = -1044.14287

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5731.815
$ws.Range("J70").Value = 5780.7075
$ws.Range("L70").Value = 5780.7075
$ws.Range("N70").Value = -6320.7075
$ws.Range("H73").Value = 5731.815
$ws.Range("J73").Value = 5780.7075
$ws.Range("L73").Value = 5780.7075
$ws.Range("N73").Value = -7652.7075
$ws.Range("H107").Value = 275.07693
$ws.Range("I107").Value = 104.6
$ws.Range("J107").Value = 381.625
$ws.Range("K107").Value = 104.6
$ws.Range("L107").Value = 381.625
$ws.Range("M107").Value = 1815.4
$ws.Range("N107").Value = -4221.625
$ws.Range("H123").Value = 26271.428
$ws.Range("J123").Value = 26271.428
$ws.Range("L123").Value = 26271.428
$ws.Range("N123").Value = -31171.428
$ws.Range("H132").Value = 1818.3334
$ws.Range("I132").Value = 1477.625
$ws.Range("J132").Value = 2499.75
$ws.Range("K132").Value = 4432.875
$ws.Range("L132").Value = 7499.25
$ws.Range("M132").Value = -1902.875
$ws.Range("N132").Value = -12559.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 743.75
$ws.Range("J22").Value = 591.6667
$ws.Range("L22").Value = 591.6667
$ws.Range("N22").Value = -1181.6667
$ws.Range("H27").Value = 743.75
$ws.Range("J27").Value = 591.6667
$ws.Range("L27").Value = 591.6667
$ws.Range("N27").Value = -805.6667
$ws.Range("H32").Value = 468.5
$ws.Range("I32").Value = 468.5
$ws.Range("K32").Value = 468.5
$ws.Range("M32").Value = -151.5
$ws.Range("H68").Value = 998
$ws.Range("I68").Value = 999.0909
$ws.Range("J68").Value = 995
$ws.Range("K68").Value = 999.0909
$ws.Range("L68").Value = 995
$ws.Range("M68").Value = -250.0909
$ws.Range("N68").Value = -2493
$ws.Range("H71").Value = 998
$ws.Range("I71").Value = 999.0909
$ws.Range("J71").Value = 995
$ws.Range("K71").Value = 4995.4545
$ws.Range("L71").Value = 4975
$ws.Range("M71").Value = -1251.4545
$ws.Range("N71").Value = -12463
$ws.Range("H132").Value = 3970.5
$ws.Range("I132").Value = 3811.1428
$ws.Range("J132").Value = 4109.9375
$ws.Range("K132").Value = 11433.4284
$ws.Range("L132").Value = 12329.8125
$ws.Range("M132").Value = -8903.428400000001
$ws.Range("N132").Value = -17389.8125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H43").Value = 9750
$ws.Range("I43").Value = 6000
$ws.Range("J43").Value = 10500
$ws.Range("K43").Value = 6000
$ws.Range("L43").Value = 10500
$ws.Range("M43").Value = -5851
$ws.Range("N43").Value = -10798
$ws.Range("H136").Value = 5370.5107
$ws.Range("I136").Value = 2774.9565
$ws.Range("J136").Value = 7857.9165
$ws.Range("K136").Value = 8324.869499999999
$ws.Range("L136").Value = 23573.7495
$ws.Range("M136").Value = -5774.869499999999
$ws.Range("N136").Value = -28673.7495
